$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Course Title (A), Category (B), Release Time (D) for rows 2-21
$ws.Range("A2").Value = "Learn Embarcadero Borland C++ Builder in 1 hour"
$ws.Range("B2").Value = "IT & Software > Other IT & Software"
$ws.Range("D2").Value = "8 hours ago"
$ws.Range("A3").Value = "Exam MS-900: Microsoft 365 Fundamentals Mock Exams"
$ws.Range("B3").Value = "IT & Software > IT Certifications"
$ws.Range("D3").Value = "8 hours ago"
$ws.Range("A4").Value = "PCEP (30-02) Practice Exams"
$ws.Range("B4").Value = "IT & Software > IT Certifications"
$ws.Range("D4").Value = "7 hours ago"
$ws.Range("A5").Value = "Linux Mastery: CLI & Kali Commands Practice Tests 2024 pro"
$ws.Range("B5").Value = "IT & Software > IT Certifications"
$ws.Range("D5").Value = "7 hours ago"
$ws.Range("A6").Value = "(ISC)2 Certified in Cybersecurity (CC) Practice Exams"
$ws.Range("B6").Value = "IT & Software > IT Certifications"
$ws.Range("D6").Value = "7 hours ago"
$ws.Range("A7").Value = "(ISC)2 Certified in Cybersecurity (CC) Practice Exams: Set 2"
$ws.Range("B7").Value = "IT & Software > IT Certifications"
$ws.Range("D7").Value = "7 hours ago"
$ws.Range("A8").Value = "CompTIA Security+ (SY0-701) Practice Tests"
$ws.Range("B8").Value = "IT & Software > IT Certifications"
$ws.Range("D8").Value = "7 hours ago"
$ws.Range("A9").Value = "CSS And Javascript Crash Course"
$ws.Range("B9").Value = "IT & Software > IT Certifications"
$ws.Range("D9").Value = "5 hours ago"
$ws.Range("A10").Value = "ECCouncil: Certified Cybersecurity Technician"
$ws.Range("B10").Value = "IT & Software > IT Certifications"
$ws.Range("D10").Value = "29 minutes ago"
$ws.Range("A11").Value = "AZ-900 Azure Fundamentals - Microsoft Azure Fundamentals"
$ws.Range("B11").Value = "IT & Software > IT Certifications"
$ws.Range("D11").Value = "2 hours ago"
$ws.Range("A12").Value = "Salesforce Certified Platform Developer I 2023"
$ws.Range("B12").Value = "IT & Software > IT Certifications"
$ws.Range("D12").Value = "17 minutes ago"
$ws.Range("A13").Value = "The Best ChatGPT & AI Course: Make Money With AI"
$ws.Range("B13").Value = "IT & Software > Other IT & Software"
$ws.Range("D13").Value = "13 hours ago"
$ws.Range("A14").Value = "Midjourney for Beginners: Embark on Your Artistic Journey"
$ws.Range("B14").Value = "IT & Software > Other IT & Software"
$ws.Range("D14").Value = "12 hours ago"
$ws.Range("A15").Value = "Learn Azure Bicep"
$ws.Range("B15").Value = "IT & Software > Other IT & Software"
$ws.Range("D15").Value = "12 hours ago"
$ws.Range("A16").Value = "Google Forms o Formularios de Cero a Avanzado"
$ws.Range("B16").Value = "IT & Software > Other IT & Software"
$ws.Range("D16").Value = "11 hours ago"
$ws.Range("A17").Value = "18 Crucial Cyber Security Tips"
$ws.Range("B17").Value = "IT & Software > Network & Security"
$ws.Range("D17").Value = "10 hours ago"
$ws.Range("A18").Value = "Web Applications Step by Step Guide Part 4"
$ws.Range("B18").Value = "IT & Software > Other IT & Software"
$ws.Range("D18").Value = "10 hours ago"
$ws.Range("A19").Value = "Web Applications Step by Step Guide Part-2"
$ws.Range("B19").Value = "IT & Software > Other IT & Software"
$ws.Range("D19").Value = "10 hours ago"
$ws.Range("A20").Value = "Web Application: Step by Step Guide"
$ws.Range("B20").Value = "IT & Software > Other IT & Software"
$ws.Range("D20").Value = "10 hours ago"
$ws.Range("A21").Value = "Web Applications Step by Step Guide Part - 3"
$ws.Range("B21").Value = "IT & Software > Other IT & Software"
$ws.Range("D21").Value = "10 hours ago"

# Rebuild hyperlinks in column C with the updated target URLs
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.real.discount/offer/learn-embarcadero-borland-c-builder-in-1-hour-2")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.real.discount/offer/exam-ms-900-microsoft-365-fundamentals-mock-exams-35283")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.real.discount/offer/pcep-30-02-practice-exams-36925")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.real.discount/offer/linux-mastery-cli-kali-commands-practice-tests-2024-pro-38400")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.real.discount/offer/isc-2-certified-in-cybersecurity-cc-practice-exams-36420")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://www.real.discount/offer/isc-2-certified-in-cybersecurity-cc-practice-exams-set-2-36444")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.real.discount/offer/comptia-security-sy0-701-practice-tests-36461")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://www.real.discount/offer/css-and-javascript-crash-course-13022")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.real.discount/offer/eccouncil-certified-cybersecurity-technician-39209")
$ws.Hyperlinks.Add($ws.Range("C11"), "https://www.real.discount/offer/az-900-azure-fundamentals-microsoft-azure-fundamentals-15711")
$ws.Hyperlinks.Add($ws.Range("C12"), "https://www.real.discount/offer/salesforce-certified-platform-developer-i-2023-37805")
$ws.Hyperlinks.Add($ws.Range("C13"), "https://www.real.discount/offer/the-best-chatgpt-ai-course-make-money-with-ai-35563")
$ws.Hyperlinks.Add($ws.Range("C14"), "https://www.real.discount/offer/midjourney-for-beginners-embark-on-your-artistic-journey-35457")
$ws.Hyperlinks.Add($ws.Range("C15"), "https://www.real.discount/offer/learn-azure-bicep-20512")
$ws.Hyperlinks.Add($ws.Range("C16"), "https://www.real.discount/offer/google-forms-o-formularios-de-cero-a-avanzado-32061")
$ws.Hyperlinks.Add($ws.Range("C17"), "https://www.real.discount/offer/18-crucial-cyber-security-tips-29894")
$ws.Hyperlinks.Add($ws.Range("C18"), "https://www.real.discount/offer/web-applications-step-by-step-guide-part-4-27521")
$ws.Hyperlinks.Add($ws.Range("C19"), "https://www.real.discount/offer/web-applications-step-by-step-guide-part-2-27523")
$ws.Hyperlinks.Add($ws.Range("C20"), "https://www.real.discount/offer/web-application-step-by-step-guide-27524")
$ws.Hyperlinks.Add($ws.Range("C21"), "https://www.real.discount/offer/web-applications-step-by-step-guide-part-3-27522")

# Restore the Hyperlink cell style so it matches the original formatting
$ws.Range("C2:C21").Style = "Hyperlink"
